$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns W1 through AO1 (inserting "Left_vs_Right" after
# "Av_Diad_HB_prom_ratio" and appending "HB1_prom/std_betweendiads" /
# "HB2_prom/std_betweendiads" at the end of the row).
$ws.Range("W1").Value = "Diad1_HB1_Valley_prom"
$ws.Range("X1").Value = "Mean_Diad_HB_Valley_prom"
$ws.Range("Y1").Value = "Diad1_prom/std_betweendiads"
$ws.Range("Z1").Value = "Diad2_prom/std_betweendiads"
$ws.Range("AA1").Value = "Av_Diad_prom/std_betweendiads"
$ws.Range("AB1").Value = "C13_prom/HB2_prom"
$ws.Range("AC1").Value = "Av_Diad_HB_prom_ratio"
$ws.Range("AD1").Value = "Left_vs_Right"
$ws.Range("AE1").Value = "Diad2_height"
$ws.Range("AF1").Value = "HB2_height"
$ws.Range("AG1").Value = "C13_height"
$ws.Range("AH1").Value = "Diad1_height"
$ws.Range("AI1").Value = "HB1_height"
$ws.Range("AJ1").Value = "Diad1_Median_Bck"
$ws.Range("AK1").Value = "Diad2_Median_Bck"
$ws.Range("AL1").Value = "C13_HB2_abs_prom_ratio"
$ws.Range("AM1").Value = "Diad2_HB2_Valley_prom"
$ws.Range("AN1").Value = "HB1_prom/std_betweendiads"
$ws.Range("AO1").Value = "HB2_prom/std_betweendiads"

# Column AP is no longer used; remove it entirely so the sheet's used range
# shrinks back to B1:AO1.
$ws.Range("AP1").EntireColumn.Delete()
